$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 758, shifting existing rows 758:799 down to 759:800
$ws.Rows.Item(758).Insert()

# Populate the newly inserted row 758 with its data.
# Column A holds a date-like string ("2026/02/05"); format the cell as Text first
# so Excel stores it as a literal string instead of auto-converting it to a date
# serial number, then clear the formatting again so the cell keeps the sheet's
# default (unstyled) look, matching the rest of the column.
$ws.Range("A758").NumberFormat = "@"
$ws.Range("A758").Value = "2026/02/05"
$ws.Range("A758").ClearFormats()

$ws.Range("B758").Value = "木"
$ws.Range("C758").Value = 10
$ws.Range("D758").Value = 39
